$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 36.35426266666666
$ws.Range("H2").Value = 109.062788
$ws.Range("I2").Value = 0.4094848412143908
$ws.Range("J2").Value = 0.4094848412143908
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 1923.515077271096
$ws.Range("R2").Value = 17311.63569543987
$ws.Range("S2").Value = 0.1704036534764422
$ws.Range("T2").Value = 0.1704036534764422

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 36.35426266666666
$ws.Range("H3").Value = 109.062788
$ws.Range("I3").Value = 0.4094848412143908
$ws.Range("J3").Value = 0.4094848412143908
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 1718.104574291988
$ws.Range("R3").Value = 15462.9411686279
$ws.Range("S3").Value = 0.1522063954545648
$ws.Range("T3").Value = 0.1522063954545648

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 36.35426266666666
$ws.Range("H4").Value = 109.062788
$ws.Range("I4").Value = 0.4094848412143908
$ws.Range("J4").Value = 0.4094848412143908
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 980.6419603261928
$ws.Range("R4").Value = 8825.777642935736
$ws.Range("S4").Value = 0.08687479228338386
$ws.Range("T4").Value = 0.08687479228338386

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 45.11545066666667
$ws.Range("H5").Value = 135.346352
$ws.Range("I5").Value = 0.5081685556916724
$ws.Range("J5").Value = 0.5081685556916724
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 2387.072194831853
$ws.Range("R5").Value = 21483.64975348668
$ws.Range("S5").Value = 0.2114700466442191
$ws.Range("T5").Value = 0.2114700466442191

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.11545066666667
$ws.Range("H6").Value = 135.346352
$ws.Range("I6").Value = 0.5081685556916724
$ws.Range("J6").Value = 0.5081685556916724
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 2132.158830241289
$ws.Range("R6").Value = 19189.4294721716
$ws.Range("S6").Value = 0.1888873441952055
$ws.Range("T6").Value = 0.1888873441952055

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.11545066666667
$ws.Range("H7").Value = 135.346352
$ws.Range("I7").Value = 0.5081685556916724
$ws.Range("J7").Value = 0.5081685556916724
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 1216.971566399705
$ws.Range("R7").Value = 10952.74409759735
$ws.Range("S7").Value = 0.1078111648522478
$ws.Range("T7").Value = 0.1078111648522478

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.310771333333334
$ws.Range("H8").Value = 21.932314
$ws.Range("I8").Value = 0.0823466030939367
$ws.Range("J8").Value = 0.0823466030939367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 386.8151312842282
$ws.Range("R8").Value = 3481.336181558054
$ws.Range("S8").Value = 0.03426784243579508
$ws.Range("T8").Value = 0.03426784243579508

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.310771333333334
$ws.Range("H9").Value = 21.932314
$ws.Range("I9").Value = 0.0823466030939367
$ws.Range("J9").Value = 0.0823466030939367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 345.5074796749944
$ws.Range("R9").Value = 3109.56731707495
$ws.Range("S9").Value = 0.0306084093313082
$ws.Range("T9").Value = 0.0306084093313082

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.310771333333334
$ws.Range("H10").Value = 21.932314
$ws.Range("I10").Value = 0.0823466030939367
$ws.Range("J10").Value = 0.0823466030939367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 197.2051860204565
$ws.Range("R10").Value = 1774.846674184108
$ws.Range("S10").Value = 0.01747035132683342
$ws.Range("T10").Value = 0.01747035132683342
